# Add Filter tahun dan bulan
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Phase 1: shift the existing formatting of row 1 one column to the
#     left (B1:M1 -> A1:L1), since the new "tahun" (year) filter column is
#     being inserted at the front of the table ---
$ws.Range("B1:M1").Copy()
$ws.Range("A1:L1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

# M1 (the new "bulan" / month column) keeps the default (unstyled) look,
# so copy the format from a cell that already uses the default style
# (C1, after the phase-1 shift, originally D1's default style)
$ws.Range("C1").Copy()
$ws.Range("M1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

# --- Phase 2: write the refreshed values for row 1 ---
$ws.Range("A1").Value = "Active"
$ws.Range("B1").Value = 1223
$ws.Range("C1").Value = 1
$ws.Range("D1").Value = 2121233
$ws.Range("E1").Value = 122132
$ws.Range("F1").Value = 2121
$ws.Range("G1").Value = 1
$ws.Range("H1").ClearContents()
$ws.Range("I1").ClearContents()
$ws.Range("J1").ClearContents()
$ws.Range("K1").ClearContents()
$ws.Range("L1").Value = "we"
$ws.Range("M1").Value = "telkomsel"

# --- Phase 3: build row 2 (another filter row) from row 1's formatting,
#     then adjust its values (G2 differs from G1) ---
$ws.Range("A1:G1").Copy()
$ws.Range("A2").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

$ws.Range("L1:M1").Copy()
$ws.Range("L2").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

$ws.Range("A2").Value = "Active"
$ws.Range("B2").Value = 1223
$ws.Range("C2").Value = 1
$ws.Range("D2").Value = 2121233
$ws.Range("E2").Value = 122132
$ws.Range("F2").Value = 2121
$ws.Range("G2").Value = 2
$ws.Range("L2").Value = "we"
$ws.Range("M2").Value = "telkomsel"

# --- Phase 4: match the saved selection / active cell ---
$ws.Range("C2").Select()
